$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (quarter labels), shifted forward by one quarter
$ws.Range("C1").Value = "2024 Q4"
$ws.Range("D1").Value = "2025 Q1"
$ws.Range("E1").Value = "2025 Q2"
$ws.Range("F1").Value = "2025 Q3"
$ws.Range("G1").Value = "2025 Q4"
$ws.Range("H1").Value = "2026 Q1"
$ws.Range("I1").Value = "2026 Q2"
$ws.Range("J1").Value = "2026 Q3"
$ws.Range("K1").Value = "2026 Q4"

# Update forecast data cells with refreshed simulation results
# Row 2
$ws.Range("C2").Value = 463.5413257
$ws.Range("D2").Value = 447.088159569206
$ws.Range("E2").Value = 442.796289684976
$ws.Range("F2").Value = 445.709739995982
$ws.Range("G2").Value = 444.7661988
$ws.Range("H2").Value = 447.863487871975
$ws.Range("I2").Value = 448.188295792375
$ws.Range("J2").Value = 452.346923275822
$ws.Range("K2").Value = 454.536674

# Row 3
$ws.Range("C3").Value = 79.275
$ws.Range("D3").Value = 79.275
$ws.Range("E3").Value = 75.203
$ws.Range("F3").Value = 75.203
$ws.Range("G3").Value = 75.203
$ws.Range("H3").Value = 75.203
$ws.Range("I3").Value = 75.203
$ws.Range("J3").Value = 75.203
$ws.Range("K3").Value = 75.203

# Row 4
$ws.Range("C4").Value = 1922
$ws.Range("D4").Value = 1915.80213881562
$ws.Range("E4").Value = 1915.49537974805
$ws.Range("F4").Value = 1921.04530428972
$ws.Range("G4").Value = 1926.83504223656
$ws.Range("H4").Value = 1936.90542286818
$ws.Range("I4").Value = 1947.35251893945
$ws.Range("J4").Value = 1958.55020912426
$ws.Range("K4").Value = 1968.55010329008

# Row 5
$ws.Range("C5").Value = 3174.7
$ws.Range("D5").Value = 3204.80183566645
$ws.Range("E5").Value = 3231.06115711615
$ws.Range("F5").Value = 3260.58967717546
$ws.Range("G5").Value = 3291.22060203007
$ws.Range("H5").Value = 3322.72316165216
$ws.Range("I5").Value = 3352.0947052084
$ws.Range("J5").Value = 3382.01967066365
$ws.Range("K5").Value = 3411.2234432123

# Row 6
$ws.Range("C6").Value = 93.901
$ws.Range("D6").Value = 93.901
$ws.Range("E6").Value = 85.879
$ws.Range("F6").Value = 85.879
$ws.Range("G6").Value = 90.218
$ws.Range("H6").Value = 88.068
$ws.Range("I6").Value = 88.068
$ws.Range("J6").Value = 88.068
$ws.Range("K6").Value = 88.068

# Row 7
$ws.Range("C7").Value = -0.901
$ws.Range("D7").Value = -0.901
$ws.Range("E7").Value = -0.901
$ws.Range("F7").Value = -0.901
$ws.Range("G7").Value = -2.15
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

# Row 8
$ws.Range("C8").Value = 0.399999999999999
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# Row 9
$ws.Range("C9").Value = 35.7
$ws.Range("D9").Value = 35.7
$ws.Range("E9").Value = 35.7
$ws.Range("F9").Value = 35.7
$ws.Range("G9").Value = 35.7
$ws.Range("H9").Value = 35.7
$ws.Range("I9").Value = 35.7
$ws.Range("J9").Value = 35.7
$ws.Range("K9").Value = 35.7

# Row 10
$ws.Range("C10").Value = 655.175
$ws.Range("D10").Value = 661.271179816286
$ws.Range("E10").Value = 674.725766162213
$ws.Range("F10").Value = 689.587400692901
$ws.Range("G10").Value = 698.873554348416
$ws.Range("H10").Value = 709.015878706698
$ws.Range("I10").Value = 719.413289541699
$ws.Range("J10").Value = 729.991705539104
$ws.Range("K10").Value = 740.621831915836

# Row 11
$ws.Range("C11").Value = 957.6
$ws.Range("D11").Value = 967.107196423835
$ws.Range("E11").Value = 986.808063219515
$ws.Range("F11").Value = 1006.91025486723
$ws.Range("G11").Value = 1021.54071990454
$ws.Range("H11").Value = 1036.38376645662
$ws.Range("I11").Value = 1051.44248334533
$ws.Range("J11").Value = 1066.72000427332
$ws.Range("K11").Value = 1082.45993400371

# Row 12
$ws.Range("C12").Value = 1118.5
$ws.Range("D12").Value = 1148.744184989
$ws.Range("E12").Value = 1163.37353203953
$ws.Range("F12").Value = 1178.18918496906
$ws.Range("G12").Value = 1195.82078219596
$ws.Range("H12").Value = 1213.71623621661
$ws.Range("I12").Value = 1231.87949564704
$ws.Range("J12").Value = 1250.31456819432
$ws.Range("K12").Value = 1269.18692299933

# Row 13
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Row 14
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0

# Row 15
$ws.Range("C15").Value = 2.372
$ws.Range("D15").Value = 2.372
$ws.Range("E15").Value = 2.372
$ws.Range("F15").Value = 2.372
$ws.Range("G15").Value = 0.49
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0

# Row 16
$ws.Range("C16").Value = 1.63
$ws.Range("D16").Value = 1.63
$ws.Range("E16").Value = 1.63
$ws.Range("F16").Value = 1.63
$ws.Range("G16").Value = 1.671
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

# Row 17
$ws.Range("C17").Value = 2171.498
$ws.Range("D17").Value = 2216.251646875
$ws.Range("E17").Value = 2222.951646875
$ws.Range("F17").Value = 2229.651646875
$ws.Range("G17").Value = 2220.112646875
$ws.Range("H17").Value = 2264.21819616788
$ws.Range("I17").Value = 2270.19019616788
$ws.Range("J17").Value = 2277.89019616788
$ws.Range("K17").Value = 2285.59019616788

# Row 18
$ws.Range("C18").Value = 224.743658420665
$ws.Range("D18").Value = 228.029047409692
$ws.Range("E18").Value = 231.362463474923
$ws.Range("F18").Value = 234.744608694576
$ws.Range("G18").Value = 238.176195410118
$ws.Range("H18").Value = 241.657946376297
$ws.Range("I18").Value = 245.190594913367
$ws.Range("J18").Value = 248.77488506154
$ws.Range("K18").Value = 252.41157173769

# Row 19
$ws.Range("C19").Value = 4526.3
$ws.Range("D19").Value = 4599.01056117126
$ws.Range("E19").Value = 4622.15249842381
$ws.Range("F19").Value = 4645.50799100346
$ws.Range("G19").Value = 4669.0771571723
$ws.Range("H19").Value = 4788.29147684674
$ws.Range("I19").Value = 4821.14923693197
$ws.Range("J19").Value = 4854.37939110322
$ws.Range("K19").Value = 4887.98436710034

# Row 20
$ws.Range("C20").Value = 2515.1
$ws.Range("D20").Value = 2529.55898292532
$ws.Range("E20").Value = 2561.00840480961
$ws.Range("F20").Value = 2588.01712691122
$ws.Range("G20").Value = 2614.62275452149
$ws.Range("H20").Value = 2639.09180766924
$ws.Range("I20").Value = 2664.54486008583
$ws.Range("J20").Value = 2687.78303814039
$ws.Range("K20").Value = 2712.94914769761

# Row 21
$ws.Range("C21").Value = 509.332755158353
$ws.Range("D21").Value = 526.099889428512
$ws.Range("E21").Value = 543.418994465887
$ws.Range("F21").Value = 561.308241039732
$ws.Range("G21").Value = 557.839987860135
$ws.Range("H21").Value = 554.393164581683
$ws.Range("I21").Value = 550.967638791707
$ws.Range("J21").Value = 547.563278895698
$ws.Range("K21").Value = 544.179954112255

# Row 22
$ws.Range("C22").Value = 163.324749650848
$ws.Range("D22").Value = 163.759599441357
$ws.Range("E22").Value = 163.691956140612
$ws.Range("F22").Value = 163.832074406442
$ws.Range("G22").Value = 163.749936112679
$ws.Range("H22").Value = 163.948034350578
$ws.Range("I22").Value = 164.27175586129
$ws.Range("J22").Value = 164.107479273765
$ws.Range("K22").Value = 164.305577511663

# Row 23
$ws.Range("C23").Value = 2.649835
$ws.Range("D23").Value = 2.645442
$ws.Range("E23").Value = 2.6696
$ws.Range("F23").Value = 2.693006
$ws.Range("G23").Value = 2.71754
$ws.Range("H23").Value = 2.71754
$ws.Range("I23").Value = 2.71754
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0

# Row 24
$ws.Range("C24").Value = 71.1791672183253
$ws.Range("D24").Value = 73.9806362924226
$ws.Range("E24").Value = 88.3944428040623
$ws.Range("F24").Value = 103.083510693309
$ws.Range("G24").Value = 109.658761987664
$ws.Range("H24").Value = 134.709036355557
$ws.Range("I24").Value = 155.206844624839
$ws.Range("J24").Value = 162.820344255939
$ws.Range("K24").Value = 0
